# Update the Approved/Rejected (I) and ReasonToReject (J) columns for the
# two test cases that were re-reviewed: change the verdict from "Rejected"
# to "Approved" and clear out the now-unused "Nil" reason.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I9").Value = "Approved"
$ws.Range("J9").ClearContents()

$ws.Range("I14").Value = "Approved"
$ws.Range("J14").ClearContents()

# Update the active selection to reflect where the user clicked next.
$ws.Range("I18").Select()
